$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6142.5884
$ws.Range("I28").Value = 6494.875
$ws.Range("K28").Value = 6494.875
$ws.Range("M28").Value = -6009.875

$ws.Range("H82").Value = 559.1429000000001
$ws.Range("I82").Value = 559.1429000000001
$ws.Range("K82").Value = 1677.4287
$ws.Range("M82").Value = -1271.4287

$ws.Range("H85").Value = 559.1429000000001
$ws.Range("I85").Value = 559.1429000000001
$ws.Range("K85").Value = 1677.4287
$ws.Range("M85").Value = -273.4287000000002

$ws.Range("H98").Value = 2240
$ws.Range("I98").Value = 2240
$ws.Range("K98").Value = 2240
$ws.Range("M98").Value = -742

$ws.Range("H113").Value = 2543.652
$ws.Range("I113").Value = 2400
$ws.Range("J113").Value = 2557.3333
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 2557.3333
$ws.Range("M113").Value = 854
$ws.Range("N113").Value = -9065.3333

$ws.Range("H122").Value = 2240
$ws.Range("I122").Value = 2240
$ws.Range("K122").Value = 6720
$ws.Range("M122").Value = -4270

$ws.Range("H132").Value = 10421088
$ws.Range("I132").Value = 15153877
$ws.Range("K132").Value = 45461631
$ws.Range("M132").Value = -45459101

$ws.Range("H138").Value = 2177.27
$ws.Range("I138").Value = 1452.4546
$ws.Range("J138").Value = 2266.854
$ws.Range("K138").Value = 4357.3638
$ws.Range("L138").Value = 6800.562
$ws.Range("M138").Value = 782.6361999999999
$ws.Range("N138").Value = -17080.562

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3039.7021
$ws.Range("I32").Value = 3200.1592
$ws.Range("K32").Value = 3200.1592
$ws.Range("M32").Value = -2913.1592

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 55000
$ws.Range("J62").Value = 55000
$ws.Range("L62").Value = 55000
$ws.Range("N62").Value = -56372

$ws.Range("H65").Value = 55000
$ws.Range("J65").Value = 55000
$ws.Range("L65").Value = 165000
$ws.Range("N65").Value = -171864

$ws.Range("H105").Value = 66669400
$ws.Range("I105").Value = 111113660
$ws.Range("K105").Value = 111113660
$ws.Range("M105").Value = -111111913

$ws.Range("H107").Value = 1640.6923
$ws.Range("I107").Value = 1314.625
$ws.Range("K107").Value = 1314.625
$ws.Range("M107").Value = 605.375

$ws.Range("H134").Value = 8409.0625
$ws.Range("I134").Value = 2162.1667
$ws.Range("K134").Value = 6486.500100000001
$ws.Range("M134").Value = -3951.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 873.4915
$ws.Range("I31").Value = 744.51166
$ws.Range("J31").Value = 1220.125
$ws.Range("K31").Value = 744.51166
$ws.Range("L31").Value = 1220.125
$ws.Range("M31").Value = -449.51166
$ws.Range("N31").Value = -1810.125

$ws.Range("H34").Value = 873.4915
$ws.Range("I34").Value = 744.51166
$ws.Range("J34").Value = 1220.125
$ws.Range("K34").Value = 744.51166
$ws.Range("L34").Value = 1220.125
$ws.Range("M34").Value = -542.51166
$ws.Range("N34").Value = -1624.125

$ws.Range("H41").Value = 6534.5
$ws.Range("I41").Value = 3610.8572
$ws.Range("K41").Value = 3610.8572
$ws.Range("M41").Value = -3182.8572

$ws.Range("H58").Value = 955.45

$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 4445.25
$ws.Range("I132").Value = 4726.5557
$ws.Range("K132").Value = 14179.6671
$ws.Range("M132").Value = -11649.6671

$ws.Range("H134").Value = 10102304
$ws.Range("I134").Value = 11495416
$ws.Range("K134").Value = 34486248
$ws.Range("M134").Value = -34483713

$ws.Range("H136").Value = 955.45

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1800.1904
$ws.Range("I126").Value = 1430
$ws.Range("J126").Value = 2725.6667
$ws.Range("K126").Value = 4290
$ws.Range("L126").Value = 8177.000100000001
$ws.Range("M126").Value = -1820
$ws.Range("N126").Value = -13117.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2336.1428
$ws.Range("I7").Value = 2068.8
$ws.Range("J7").Value = 3004.5
$ws.Range("K7").Value = 2068.8
$ws.Range("L7").Value = 3004.5
$ws.Range("M7").Value = -1956.8
$ws.Range("N7").Value = -3228.5

$ws.Range("H40").Value = 2563
$ws.Range("I40").Value = 2199
$ws.Range("K40").Value = 2199
$ws.Range("M40").Value = -2063

$ws.Range("H61").Value = 2421.05
$ws.Range("I61").Value = 1922.2142
$ws.Range("K61").Value = 1922.2142
$ws.Range("M61").Value = -1720.2142

$ws.Range("H113").Value = 2421.05
$ws.Range("I113").Value = 1922.2142
$ws.Range("K113").Value = 1922.2142
$ws.Range("M113").Value = 247.7858000000001

$ws.Range("H122").Value = 31252302
$ws.Range("I122").Value = 35716490
$ws.Range("K122").Value = 107149470
$ws.Range("M122").Value = -107147020

$ws.Range("H126").Value = 2336.1428
$ws.Range("I126").Value = 2068.8
$ws.Range("J126").Value = 3004.5
$ws.Range("K126").Value = 6206.400000000001
$ws.Range("L126").Value = 9013.5
$ws.Range("M126").Value = -3736.400000000001
$ws.Range("N126").Value = -13953.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1571.4286
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 2333.3333
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 4666.6666
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -5748.6666

$ws.Range("H107").Value = 925
$ws.Range("J107").Value = 1500
$ws.Range("L107").Value = 4500
$ws.Range("N107").Value = -8340

$ws.Range("H132").Value = 1744.0588
$ws.Range("I132").Value = 1403.3334
$ws.Range("J132").Value = 4299.5
$ws.Range("K132").Value = 4210.0002
$ws.Range("L132").Value = 12898.5
$ws.Range("M132").Value = -1680.0002
$ws.Range("N132").Value = -17958.5

$ws.Range("H136").Value = 1850.0714
$ws.Range("I136").Value = 1686
$ws.Range("J136").Value = 2014.1428
$ws.Range("K136").Value = 5058
$ws.Range("L136").Value = 6042.428400000001
$ws.Range("M136").Value = -2508
$ws.Range("N136").Value = -11142.4284
